$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D12").Value = "“케라스 창시자에게 배우는 딥러닝” 텐서플로 2.4 업데이트"
$ws.Range("E12").Value = "https://tensorflow.blog/2021/03/18/%ec%bc%80%eb%9d%bc%ec%8a%a4-%ec%b0%bd%ec%8b%9c%ec%9e%90%ec%97%90%ea%b2%8c-%eb%b0%b0%ec%9a%b0%eb%8a%94-%eb%94%a5%eb%9f%ac%eb%8b%9d-%ed%85%90%ec%84%9c%ed%94%8c%eb%a1%9c-2-4-%ec%97%85%eb%8d%b0/"

$ws.Range("D32").Value = "Maximum Likelihood Estimation (Simple Error Bound)"
$ws.Range("E32").Value = "https://dodonam.tistory.com/304"

$ws.Range("D39").Value = "How to Use Poisson Distribution like You Know What You Are Doing"
$ws.Range("E39").Value = "https://a292run.tistory.com/entry/How-to-Use-Poisson-Distribution-like-You-Know-What-You-Are-Doing-1"

$ws.Range("D51").Value = "카카오 책 검색 API 이용해서 무한 스크롤 기능 가진 도서 검색 웹 페이지 만들기"
$ws.Range("E51").Value = "https://bskyvision.com/1148"
